$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers in E1 and F1
$ws.Range("E1").Value = "FirstName"
$ws.Range("F1").Value = "LastName"

# Add new data in E2 and F2
$ws.Range("E2").Value = "Yashwanth"
$ws.Range("F2").Value = "Arul"

# Set column widths for new columns E and F to match target layout
$ws.Columns.Item(5).ColumnWidth = 13.7109375
$ws.Columns.Item(6).ColumnWidth = 9.7109375

# Update selection to C4 as per diff
$ws.Range("C4").Select()
